# Daily attendance processing - 2025-10-30 12:39:14
# Normalizes the ordering of names/emails in column G ("Attendance Recorded By"
# style list) on the "Session Analysis Results" sheet: cells that start with
# "System, " followed by one or more other identifiers are rewritten so that
# "System" moves from the front to just before the final (email) entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows in column G whose text needs to be reordered.
$rows = @(2,3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,29,30,33,37,38,39,40,41,42,44,45,46,47,48,49,51,56,57,60,64,65,66,67,68,69,71,72,73,74,75,76,78,86,87,88,89,93,95,96,97,99,102,104,112,113,114,115,119,121,122,123,125,128,130,138,139,140,141,145,147,148,149,151,154,156)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $orig = "" + $cell.Text

    $parts = $orig.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -eq 2) {
        # "System, X" -> "X, System"
        $newValue = $parts[1] + ", " + $parts[0]
    }
    elseif ($parts.Length -eq 3) {
        # "System, Y, X" -> "Y, System, X"
        $newValue = $parts[1] + ", " + $parts[0] + ", " + $parts[2]
    }
    else {
        $newValue = $orig
    }

    $cell.Value = $newValue
}
